$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of data (new Mac-Addresses) at rows 31 and 32,
# following the same pattern as the existing rows.
$ws.Range("A31").Value = 10001
$ws.Range("B31").Value = 110030
$ws.Range("C31").Value = "eng"
$ws.Range("D31").Value = $true
$ws.Range("E31").Value = "superadmin"
$ws.Range("F31").Value = "now()"

$ws.Range("A32").Value = 10001
$ws.Range("B32").Value = 110031
$ws.Range("C32").Value = "eng"
$ws.Range("D32").Value = $true
$ws.Range("E32").Value = "superadmin"
$ws.Range("F32").Value = "now()"

# Update the visible window / selection to match the post-edit state.
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("D32").Select()
